# Auto-generated edit script: apply numeric updates to Leve profit tables
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 408.1
$ws.Range("I6").Value = 413.5
$ws.Range("J6").Value = 400
$ws.Range("K6").Value = 1240.5
$ws.Range("L6").Value = 1200
$ws.Range("M6").Value = -1128.5
$ws.Range("N6").Value = -1424

$ws.Range("H33").Value = 780.85187
$ws.Range("I33").Value = 597.6842
$ws.Range("K33").Value = 597.6842
$ws.Range("M33").Value = -368.6842

$ws.Range("H43").Value = 114488
$ws.Range("I43").Value = 5666.6665
$ws.Range("J43").Value = 161125.72
$ws.Range("K43").Value = 5666.6665
$ws.Range("L43").Value = 161125.72
$ws.Range("M43").Value = -5597.6665
$ws.Range("N43").Value = -161263.72

$ws.Range("H53").Value = 303.2143
$ws.Range("I53").Value = 293.83334
$ws.Range("J53").Value = 310.25
$ws.Range("K53").Value = 293.83334
$ws.Range("L53").Value = 310.25
$ws.Range("M53").Value = 343.16666
$ws.Range("N53").Value = -1584.25

$ws.Range("H88").Value = 2746.04
$ws.Range("I88").Value = 496
$ws.Range("J88").Value = 3308.55
$ws.Range("K88").Value = 496
$ws.Range("L88").Value = 3308.55
$ws.Range("M88").Value = -90
$ws.Range("N88").Value = -4120.55

$ws.Range("H91").Value = 2746.04
$ws.Range("I91").Value = 496
$ws.Range("J91").Value = 3308.55
$ws.Range("K91").Value = 496
$ws.Range("L91").Value = 3308.55
$ws.Range("M91").Value = 908
$ws.Range("N91").Value = -6116.55

$ws.Range("H129").Value = 658.2
$ws.Range("I129").Value = 346.6
$ws.Range("J129").Value = 1281.4
$ws.Range("K129").Value = 1039.8
$ws.Range("L129").Value = 3844.2
$ws.Range("M129").Value = 3960.2
$ws.Range("N129").Value = -13844.2

$ws.Range("H132").Value = 5771.0586
$ws.Range("I132").Value = 3562.9167
$ws.Range("K132").Value = 10688.7501
$ws.Range("M132").Value = -8158.750100000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 16121.215
$ws.Range("I32").Value = 7782.051
$ws.Range("J32").Value = 30592.117
$ws.Range("K32").Value = 7782.051
$ws.Range("L32").Value = 30592.117
$ws.Range("M32").Value = -7495.051
$ws.Range("N32").Value = -31166.117

$ws.Range("H61").Value = 3555.4358
$ws.Range("I61").Value = 2740.9167
$ws.Range("J61").Value = 4858.6665
$ws.Range("K61").Value = 2740.9167
$ws.Range("L61").Value = 4858.6665
$ws.Range("M61").Value = -2528.9167
$ws.Range("N61").Value = -5282.6665

$ws.Range("H88").Value = 2009.8611
$ws.Range("I88").Value = 2120.2
$ws.Range("J88").Value = 1759.091
$ws.Range("K88").Value = 2120.2
$ws.Range("L88").Value = 1759.091
$ws.Range("M88").Value = -1714.2
$ws.Range("N88").Value = -2571.091

$ws.Range("H91").Value = 2009.8611
$ws.Range("I91").Value = 2120.2
$ws.Range("J91").Value = 1759.091
$ws.Range("K91").Value = 2120.2
$ws.Range("L91").Value = 1759.091
$ws.Range("M91").Value = -716.1999999999998
$ws.Range("N91").Value = -4567.091

$ws.Range("H97").Value = 822.14703
$ws.Range("I97").Value = 822.5172
$ws.Range("J97").Value = 820
$ws.Range("K97").Value = 822.5172
$ws.Range("L97").Value = 820
$ws.Range("M97").Value = -326.5172
$ws.Range("N97").Value = -1812

$ws.Range("H122").Value = 2387.2812
$ws.Range("I122").Value = 2300.88
$ws.Range("J122").Value = 2695.8572
$ws.Range("K122").Value = 6902.64
$ws.Range("L122").Value = 8087.571599999999
$ws.Range("M122").Value = -4452.64
$ws.Range("N122").Value = -12987.5716

$ws.Range("H136").Value = 3555.4358
$ws.Range("I136").Value = 2740.9167
$ws.Range("J136").Value = 4858.6665
$ws.Range("K136").Value = 8222.750100000001
$ws.Range("L136").Value = 14575.9995
$ws.Range("M136").Value = -5672.750100000001
$ws.Range("N136").Value = -19675.9995

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 400.61905
$ws.Range("I80").Value = 199.8
$ws.Range("K80").Value = 199.8
$ws.Range("M80").Value = 798.2

$ws.Range("H83").Value = 400.61905
$ws.Range("I83").Value = 199.8
$ws.Range("K83").Value = 999
$ws.Range("M83").Value = 3993

$ws.Range("H86").Value = 3670.3958
$ws.Range("I86").Value = 3246.3103
$ws.Range("J86").Value = 4317.684
$ws.Range("K86").Value = 3246.3103
$ws.Range("L86").Value = 4317.684
$ws.Range("M86").Value = -2123.3103
$ws.Range("N86").Value = -6563.684

$ws.Range("H89").Value = 3670.3958
$ws.Range("I89").Value = 3246.3103
$ws.Range("J89").Value = 4317.684
$ws.Range("K89").Value = 16231.5515
$ws.Range("L89").Value = 21588.42
$ws.Range("M89").Value = -10615.5515
$ws.Range("N89").Value = -32820.42

$ws.Range("H107").Value = 1966.3334
$ws.Range("I107").Value = 2122.4167
$ws.Range("J107").Value = 1342
$ws.Range("K107").Value = 2122.4167
$ws.Range("L107").Value = 1342
$ws.Range("M107").Value = -202.4167000000002
$ws.Range("N107").Value = -5182

$ws.Range("H134").Value = 24020.318
$ws.Range("I134").Value = 31204.697
$ws.Range("J134").Value = 7085.7144
$ws.Range("K134").Value = 93614.091
$ws.Range("L134").Value = 21257.1432
$ws.Range("M134").Value = -91079.091
$ws.Range("N134").Value = -26327.1432

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H36").Value = 5000
$ws.Range("I36").Value = 5000
$ws.Range("K36").Value = 5000
$ws.Range("M36").Value = -4612

$ws.Range("H40").Value = 5000
$ws.Range("I40").Value = 5000
$ws.Range("K40").Value = 5000
$ws.Range("M40").Value = -4840

$ws.Range("H94").Value = 5854.647
$ws.Range("I94").Value = 1250
$ws.Range("J94").Value = 6468.6
$ws.Range("K94").Value = 1250
$ws.Range("L94").Value = 6468.6
$ws.Range("M94").Value = -799
$ws.Range("N94").Value = -7370.6

$ws.Range("H132").Value = 16131852
$ws.Range("I132").Value = 38462780
$ws.Range("K132").Value = 115388340
$ws.Range("M132").Value = -115385810

$ws.Range("H134").Value = 12822238
$ws.Range("I134").Value = 29412870
$ws.Range("J134").Value = 2202.6365
$ws.Range("K134").Value = 88238610
$ws.Range("L134").Value = 6607.9095
$ws.Range("M134").Value = -88236075
$ws.Range("N134").Value = -11677.9095

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 492.14816
$ws.Range("I23").Value = 75.8
$ws.Range("J23").Value = 586.7727
$ws.Range("K23").Value = 227.4
$ws.Range("L23").Value = 1760.3181
$ws.Range("M23").Value = 7.600000000000023
$ws.Range("N23").Value = -2230.3181

$ws.Range("H75").Value = 1204.2858
$ws.Range("I75").Value = 400
$ws.Range("J75").Value = 1807.5
$ws.Range("K75").Value = 1200
$ws.Range("L75").Value = 5422.5
$ws.Range("M75").Value = -202
$ws.Range("N75").Value = -7418.5

$ws.Range("H78").Value = 1204.2858
$ws.Range("I78").Value = 400
$ws.Range("J78").Value = 1807.5
$ws.Range("K78").Value = 3600
$ws.Range("L78").Value = 16267.5
$ws.Range("M78").Value = 1392
$ws.Range("N78").Value = -26251.5

$ws.Range("H86").Value = 850
$ws.Range("I86").Value = 750
$ws.Range("J86").Value = 1000
$ws.Range("K86").Value = 2250
$ws.Range("L86").Value = 3000
$ws.Range("M86").Value = -1064
$ws.Range("N86").Value = -5372

$ws.Range("H89").Value = 850
$ws.Range("I89").Value = 750
$ws.Range("J89").Value = 1000
$ws.Range("K89").Value = 6750
$ws.Range("L89").Value = 9000
$ws.Range("M89").Value = -822
$ws.Range("N89").Value = -20856

$ws.Range("H114").Value = 3706222
$ws.Range("I114").Value = 3482.2666
$ws.Range("J114").Value = 7408961.5
$ws.Range("K114").Value = 10446.7998
$ws.Range("L114").Value = 22226884.5
$ws.Range("M114").Value = -7192.799800000001
$ws.Range("N114").Value = -22233392.5

$ws.Range("H119").Value = 1973
$ws.Range("I119").Value = 859.5
$ws.Range("J119").Value = 4200
$ws.Range("K119").Value = 2578.5
$ws.Range("L119").Value = 12600
$ws.Range("M119").Value = 2259.5
$ws.Range("N119").Value = -22276

$ws.Range("H131").Value = 1925.6936
$ws.Range("I131").Value = 4938.143
$ws.Range("J131").Value = 1542.2909
$ws.Range("K131").Value = 14814.429
$ws.Range("L131").Value = 4626.8727
$ws.Range("M131").Value = -9774.429
$ws.Range("N131").Value = -14706.8727

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 1429782.4
$ws.Range("I11").Value = 5001000
$ws.Range("J11").Value = 1295.4
$ws.Range("K11").Value = 5001000
$ws.Range("L11").Value = 1295.4
$ws.Range("M11").Value = -5000861
$ws.Range("N11").Value = -1573.4

$ws.Range("H102").Value = 2529.102
$ws.Range("I102").Value = 2602.6047
$ws.Range("K102").Value = 2602.6047
$ws.Range("M102").Value = -980.6046999999999

$ws.Range("H113").Value = 4897.64
$ws.Range("I113").Value = 6270.8887
$ws.Range("J113").Value = 1366.4286
$ws.Range("K113").Value = 6270.8887
$ws.Range("L113").Value = 1366.4286
$ws.Range("M113").Value = -4100.8887
$ws.Range("N113").Value = -5706.4286

$ws.Range("H126").Value = 3118
$ws.Range("J126").Value = 2162.375
$ws.Range("L126").Value = 6487.125
$ws.Range("N126").Value = -11427.125

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 2647.5
$ws.Range("I82").Value = 2549
$ws.Range("J82").Value = 2828.0833
$ws.Range("K82").Value = 2549
$ws.Range("L82").Value = 2828.0833
$ws.Range("M82").Value = -2188
$ws.Range("N82").Value = -3550.0833

$ws.Range("H85").Value = 2647.5
$ws.Range("I85").Value = 2549
$ws.Range("J85").Value = 2828.0833
$ws.Range("K85").Value = 2549
$ws.Range("L85").Value = 2828.0833
$ws.Range("M85").Value = -1301
$ws.Range("N85").Value = -5324.0833

$ws.Range("H136").Value = 4495.7236
$ws.Range("I136").Value = 2620.2173
$ws.Range("J136").Value = 6293.0835
$ws.Range("K136").Value = 7860.651899999999
$ws.Range("L136").Value = 18879.2505
$ws.Range("M136").Value = -5310.651899999999
$ws.Range("N136").Value = -23979.2505

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 418.43478
$ws.Range("I107").Value = 414.2
$ws.Range("J107").Value = 426.375
$ws.Range("K107").Value = 1242.6
$ws.Range("L107").Value = 1279.125
$ws.Range("M107").Value = 677.4000000000001
$ws.Range("N107").Value = -5119.125
